# Clean up form formatting
# - survey sheet: mark the `source` and `name` fields (inputs group) as
#   hidden via the appearance column (F6, F9)
# - survey sheet: switch the deprecated `string` type to `text` for the
#   contact/_id field (A8)
# - leave the active selection on the cell that was last touched (F16)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("F6").Value = "hidden"
$ws.Range("A8").Value = "text"
$ws.Range("F9").Value = "hidden"

$ws.Activate()
$ws.Range("F16").Select()
